$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.417.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.76%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.653.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.87%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9980"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.96%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.35"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.61%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3285"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.118"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06922"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9999"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.948"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.58%  "

$ws.Range("E14").Value = "  -8.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.653.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.564"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001045"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06490"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9983"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "76.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.911"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.05%  "

$ws.Range("E23").Value = "  -9.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.425.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.435"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.342"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -16.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "145.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.836.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.164"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.81%  "

$ws.Range("E32").Value = "  -4.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.590"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -18.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08338"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.673"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.180"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.99%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06053"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02214"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.254"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.206"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2042"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.71%  "

$ws.Range("E43").Value = "  -0.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5849"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.720"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.53%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5568"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.98%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.931"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06895"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.63%  "
